$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 479, pushing the existing rows 479:583 down to 480:584
$ws.Rows.Item(479).Insert()

# Fill in the new weekly data row (the constant columns mirror the rest of the table)
$ws.Cells.Item(479, 1).Value = 7
$ws.Cells.Item(479, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(479, 3).Value = "Ñuble"
$ws.Cells.Item(479, 4).Value = 45209
$ws.Cells.Item(479, 5).Value = 16
$ws.Cells.Item(479, 6).Value = "Fruta"
$ws.Cells.Item(479, 7).Value = 100101
$ws.Cells.Item(479, 8).Value = "Berries"
$ws.Cells.Item(479, 9).Value = 100112025
$ws.Cells.Item(479, 10).Value = "Frutilla"
$ws.Cells.Item(479, 11).Value = "Sin especificar"
$ws.Cells.Item(479, 12).Value = "Primera"
$ws.Cells.Item(479, 13).Value = 160
$ws.Cells.Item(479, 14).Value = 11000
$ws.Cells.Item(479, 15).Value = 12000
$ws.Cells.Item(479, 16).Value = 11500
$ws.Cells.Item(479, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(479, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(479, 19).Value = 1643
$ws.Cells.Item(479, 20).Value = 7
